$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (includes swapping the Monero / ImmutableX rows 35-36)
$ws.Range("D2").Value = "68.036.03"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "3.246.15"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'581.91"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("D6").Value = "'184.87"
$ws.Range("E6").Value = "  +0.99%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +0.69%  "
$ws.Range("E9").Value = "  -3.37%  "
$ws.Range("D10").Value = "'6.61"
$ws.Range("E10").Value = "  -1.01%  "
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("D12").Value = "3.811.99"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("D14").Value = "'27.84"
$ws.Range("E14").Value = "  -2.90%  "
$ws.Range("D15").Value = "68.039.55"
$ws.Range("E15").Value = "  +0.32%  "
$ws.Range("E16").Value = "  -0.95%  "
$ws.Range("D17").Value = "3.238.80"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").Value = "'5.80"
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("D20").Value = "'396.64"
$ws.Range("E20").Value = "  +4.39%  "
$ws.Range("D21").Value = "'7.59"
$ws.Range("E21").Value = "  -0.56%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").Value = "'71.34"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  +0.59%  "
$ws.Range("E25").Value = "  -0.77%  "
$ws.Range("E26").Value = "  +2.47%  "
$ws.Range("D27").Value = "'9.62"
$ws.Range("E27").Value = "  -3.13%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("E29").Value = "  -1.15%  "
$ws.Range("D30").Value = "'5.60"
$ws.Range("E30").Value = "  -1.19%  "
$ws.Range("D31").Value = "'22.81"
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("D32").Value = "'7.02"
$ws.Range("E32").Value = "  -0.83%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").Value = "'161.84"
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'1.49"
$ws.Range("E36").Value = "  -4.48%  "
$ws.Range("E37").Value = "  +3.22%  "
$ws.Range("D38").Value = "'26.69"
$ws.Range("E38").Value = "  +0.82%  "
$ws.Range("D39").Value = "'0.812"
$ws.Range("E39").Value = "  -2.98%  "
$ws.Range("D40").Value = "'4.58"
$ws.Range("E40").Value = "  +0.26%  "
$ws.Range("D41").Value = "'6.49"
$ws.Range("E41").Value = "  -3.16%  "
$ws.Range("E42").Value = "  -3.78%  "
$ws.Range("D43").Value = "'41.20"
$ws.Range("D44").Value = "'25.28"
$ws.Range("E44").Value = "  -0.62%  "
$ws.Range("D45").Value = "'0.0684"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("D46").Value = "2.614.22"
$ws.Range("E46").Value = "  -0.42%  "
$ws.Range("D47").Value = "'335.87"
$ws.Range("E47").Value = "  -2.85%  "
$ws.Range("D48").Value = "'0.0279"
$ws.Range("E48").Value = "  -1.46%  "
$ws.Range("D49").Value = "'6.31"
$ws.Range("E49").Value = "  +2.34%  "
$ws.Range("E50").Value = "  -1.13%  "
$ws.Range("D51").Value = "'31.11"
$ws.Range("E51").Value = "  +2.32%  "
